# Auto-generated edit script: applies the "Updated cryptos list" diff
# by writing the new Price (column D) and Volume(1h) (column E) text values,
# plus the three swapped Coin/Link/Price/Volume row pairs (15/16, 34/35, 44/45).
#
# Cells in column D frequently hold strings that *look* like plain numbers
# (e.g. "15.21", "0.999"). Assigning such a string straight to .Value lets
# Excel's normal type inference kick in and silently convert it to a float
# (which then round-trips with binary noise like 15.210000000000001, and
# loses values such as the leading zeros in "0.0800"). To keep these as the
# literal text found in the source diff, a leading apostrophe is used for
# any replacement value that is a bare number -- exactly what typing
# '15.21 into a General-formatted cell does in the Excel UI -- which stores
# the text untouched while only marking the cell quotePrefix.
# Values that already aren't pure numbers (URLs, coin names, the percentage
# strings with a trailing "%", or the "42.499.61"-style multi-dot numbers)
# are assigned as plain text and need no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = '42.533.19'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '2.523.01'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''314.53'
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''98.91'
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '''0.516'
$ws.Range("E9").Value = '  -3.13%  '
$ws.Range("D10").Value = '''35.19'
$ws.Range("E10").Value = '  -3.36%  '
$ws.Range("D11").Value = '''0.0801'
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").Value = '''7.19'
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("D14").Value = '2.910.88'
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''15.21'
$ws.Range("E15").Value = '  -5.78%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.481.85'
$ws.Range("E16").Value = '  -3.97%  '
$ws.Range("D17").Value = '''0.809'
$ws.Range("E17").Value = '  -4.84%  '
$ws.Range("D18").Value = '42.525.46'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").Value = '''6.59'
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("D20").Value = '0.0₃0939'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").Value = '''12.14'
$ws.Range("E21").Value = '  -2.37%  '
$ws.Range("D22").Value = '''69.04'
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("D23").Value = '''241.50'
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("E24").Value = '  -2.66%  '
$ws.Range("E25").Value = '  -3.89%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '''25.51'
$ws.Range("E27").Value = '  -3.90%  '
$ws.Range("E28").Value = '  -4.97%  '
$ws.Range("D29").Value = '''9.99'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").Value = '''37.69'
$ws.Range("E30").Value = '  -7.71%  '
$ws.Range("D31").Value = '''5.93'
$ws.Range("E31").Value = '  +3.56%  '
$ws.Range("D32").Value = '''156.19'
$ws.Range("E32").Value = '  -1.10%  '
$ws.Range("D33").Value = '''2.70'
$ws.Range("E33").Value = '  -2.48%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0782'
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.63'
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("D37").Value = '''1.97'
$ws.Range("E37").Value = '  -5.52%  '
$ws.Range("D38").Value = '''17.51'
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("E40").Value = '  -1.00%  '
$ws.Range("D41").Value = '''4.22'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").Value = '''22.11'
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.015.63'
$ws.Range("E44").Value = '  +2.26%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '''0.0295'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '''3.21'
$ws.Range("E46").Value = '  -3.66%  '
$ws.Range("D47").Value = '''8.94'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("D48").Value = '2.765.56'
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("D49").Value = '''78.98'
$ws.Range("E49").Value = '  -3.12%  '
$ws.Range("D50").Value = '''0.188'
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("D51").Value = '''71.43'
$ws.Range("E51").Value = '  -2.97%  '
